{"js": "// The worksheet-style table holds 5 \"data\" rows (each 5 cells of\n// \"A\u00f7B=Q, R\" division answers) interleaved with blank spacer rows.\n// Row indices below are 0-based indices into Table.getCell(row, col),\n// i.e. they count every table row (blank spacers included), matching\n// Word's own TableRow numbering.\nconst replacements = [\n  { row: 0, col: 0, text: \"42\u00f77=6, 0\" },\n  { row: 0, col: 1, text: \"27\u00f79=3, 0\" },\n  { row: 0, col: 2, text: \"91\u00f72=45, 1\" },\n  { row: 0, col: 3, text: \"98\u00f76=16, 2\" },\n  { row: 0, col: 4, text: \"29\u00f77=4, 1\" },\n\n  { row: 4, col: 0, text: \"90\u00f78=11, 2\" },\n  { row: 4, col: 1, text: \"23\u00f75=4, 3\" },\n  { row: 4, col: 2, text: \"89\u00f72=44, 1\" },\n  { row: 4, col: 3, text: \"25\u00f79=2, 7\" },\n  { row: 4, col: 4, text: \"55\u00f77=7, 6\" },\n\n  { row: 8, col: 0, text: \"41\u00f77=5, 6\" },\n  { row: 8, col: 1, text: \"77\u00f78=9, 5\" },\n  { row: 8, col: 2, text: \"40\u00f72=20, 0\" },\n  { row: 8, col: 3, text: \"50\u00f78=6, 2\" },\n  { row: 8, col: 4, text: \"41\u00f74=10, 1\" },\n\n  { row: 12, col: 0, text: \"20\u00f79=2, 2\" },\n  { row: 12, col: 1, text: \"32\u00f75=6, 2\" },\n  { row: 12, col: 2, text: \"56\u00f78=7, 0\" },\n  { row: 12, col: 3, text: \"34\u00f72=17, 0\" },\n  { row: 12, col: 4, text: \"46\u00f73=15, 1\" },\n\n  { row: 16, col: 0, text: \"33\u00f76=5, 3\" },\n  { row: 16, col: 1, text: \"77\u00f76=12, 5\" },\n  { row: 16, col: 2, text: \"66\u00f72=33, 0\" },\n  { row: 16, col: 3, text: \"28\u00f78=3, 4\" },\n  { row: 16, col: 4, text: \"31\u00f72=15, 1\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Setting Cell.value directly (rather than a global Find & Replace)\n// avoids any ambiguity from the fact that a few of the new answers\n// are textually identical to other cells' old answers elsewhere in\n// the table (e.g. \"20\u00f79=2, 2\" is both an old value and a new value).\nfor (const { row, col, text } of replacements) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# The table holds 5 \"data\" rows (each with 5 cells of \"A\u00f7B=Q, R\"\n# division answers) interleaved with blank spacer rows. Rows/Cells are\n# 1-based in the Word object model and Rows.Item() counts every row,\n# including the blank spacers, so the data rows are 1, 5, 9, 13, 17.\n$replacements = @(\n    @{ Row = 1;  Col = 1; Text = \"42\u00f77=6, 0\" },\n    @{ Row = 1;  Col = 2; Text = \"27\u00f79=3, 0\" },\n    @{ Row = 1;  Col = 3; Text = \"91\u00f72=45, 1\" },\n    @{ Row = 1;  Col = 4; Text = \"98\u00f76=16, 2\" },\n    @{ Row = 1;  Col = 5; Text = \"29\u00f77=4, 1\" },\n\n    @{ Row = 5;  Col = 1; Text = \"90\u00f78=11, 2\" },\n    @{ Row = 5;  Col = 2; Text = \"23\u00f75=4, 3\" },\n    @{ Row = 5;  Col = 3; Text = \"89\u00f72=44, 1\" },\n    @{ Row = 5;  Col = 4; Text = \"25\u00f79=2, 7\" },\n    @{ Row = 5;  Col = 5; Text = \"55\u00f77=7, 6\" },\n\n    @{ Row = 9;  Col = 1; Text = \"41\u00f77=5, 6\" },\n    @{ Row = 9;  Col = 2; Text = \"77\u00f78=9, 5\" },\n    @{ Row = 9;  Col = 3; Text = \"40\u00f72=20, 0\" },\n    @{ Row = 9;  Col = 4; Text = \"50\u00f78=6, 2\" },\n    @{ Row = 9;  Col = 5; Text = \"41\u00f74=10, 1\" },\n\n    @{ Row = 13; Col = 1; Text = \"20\u00f79=2, 2\" },\n    @{ Row = 13; Col = 2; Text = \"32\u00f75=6, 2\" },\n    @{ Row = 13; Col = 3; Text = \"56\u00f78=7, 0\" },\n    @{ Row = 13; Col = 4; Text = \"34\u00f72=17, 0\" },\n    @{ Row = 13; Col = 5; Text = \"46\u00f73=15, 1\" },\n\n    @{ Row = 17; Col = 1; Text = \"33\u00f76=5, 3\" },\n    @{ Row = 17; Col = 2; Text = \"77\u00f76=12, 5\" },\n    @{ Row = 17; Col = 3; Text = \"66\u00f72=33, 0\" },\n    @{ Row = 17; Col = 4; Text = \"28\u00f78=3, 4\" },\n    @{ Row = 17; Col = 5; Text = \"31\u00f72=15, 1\" }\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Setting Cell.Range.Text directly (rather than a document-wide\n# Find & Replace) avoids any ambiguity from the fact that a few of the\n# new answers are textually identical to other cells' old answers\n# elsewhere in the table (e.g. \"20\u00f79=2, 2\" is both an old value and a\n# new value), which would otherwise make Find latch onto text that was\n# just written instead of the intended original cell.\nforeach ($item in $replacements) {\n    $cell = $tbl.Rows.Item($item.Row).Cells.Item($item.Col)\n    $cell.Range.Text = $item.Text\n}\n"}
